# Update impf_daten.xlsx: append newest BAG vaccination-data rows
# (Stand 2021-02-03, serial 44230) for all 26 Swiss cantons to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number-format (built-in short date) of the last existing
# "Stand" cell (J154) down onto the new J155:J180 range before writing
# values, so the new date cells reuse the existing date style instead of
# Excel minting a brand-new (duplicate) cell format.
$ws.Range("J154").Copy()
$ws.Range("J155:J180").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @(155, 'Genf',           'GE', 20238, 4.01, 44230),
    @(156, 'Waadt',          'VD', 34784, 4.32, 44230),
    @(157, 'Wallis',         'VS', 16735, 4.84, 44230),
    @(158, 'Freiburg',       'FR', 11769, 3.66, 44230),
    @(159, 'Neuenburg',      'NE', 5657,  3.21, 44230),
    @(160, 'Jura',           'JU', 3411,  4.64, 44230),
    @(161, 'Bern',           'BE', 31375, 3.02, 44230),
    @(162, 'Solothurn',      'SO', 16509, 6,    44230),
    @(163, 'Basel-Stadt',    'BS', 14595, 7.45, 44230),
    @(164, 'Basel-Landsch.', 'BL', 13740, 4.75, 44230),
    @(165, 'Aargau',         'AG', 24812, 3.62, 44230),
    @(166, 'Zürich',         'ZH', 54243, 3.52, 44230),
    @(167, 'Schaffhausen',   'SH', 5728,  6.96, 44230),
    @(168, 'Thurgau',        'TG', 10919, 3.91, 44230),
    @(169, 'Appenzell-A.',   'AR', 3562,  6.42, 44230),
    @(170, 'Appenzell-I.',   'AI', 1395,  8.65, 44230),
    @(171, 'St. Gallen',     'SG', 23971, 4.69, 44230),
    @(172, 'Glarus',         'GL', 2635,  6.49, 44230),
    @(173, 'Schwyz',         'SZ', 5761,  3.59, 44230),
    @(174, 'Zug',            'ZG', 7636,  5.98, 44230),
    @(175, 'Luzern',         'LU', 21035, 5.09, 44230),
    @(176, 'Nidwalden',      'NW', 3737,  8.67, 44230),
    @(177, 'Obwalden',       'OW', 2156,  5.68, 44230),
    @(178, 'Uri',            'UR', 2699,  7.35, 44230),
    @(179, 'Graubünden',     'GR', 8600,  4.32, 44230),
    @(180, 'Tessin',         'TI', 20550, 5.85, 44230)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: Kanton
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: Abk.
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: Total_Impfungen
    $ws.Cells.Item($r, 8).Value = $row[4]   # H: Impfungen pro 100 Einwohner
    $ws.Cells.Item($r, 10).Value = $row[5]  # J: Stand
}

# Match the author's final selection/scroll position from the diff.
$null = $ws.Range("L166").Select()
